$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.286.91"
$ws.Range("E2").Value = "  +5.83%  "
$ws.Range("D3").Value = "1.786.58"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("E4").Value = "  +0.18%  "
$c = $ws.Range("D5")
$s = $c.Style
$c.Value = "'246.68"
$c.Style = $s
$ws.Range("E5").Value = "  +1.84%  "
$c = $ws.Range("D6")
$s = $c.Style
$c.Value = "'0.9994"
$c.Style = $s
$ws.Range("E6").Value = "  +0.09%  "
$c = $ws.Range("D7")
$s = $c.Style
$c.Value = "'0.4915"
$c.Style = $s
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +2.56%  "
$c = $ws.Range("D9")
$s = $c.Style
$c.Value = "'0.06287"
$c.Style = $s
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "1.785.79"
$ws.Range("E10").Value = "  +3.11%  "
$c = $ws.Range("D11")
$s = $c.Style
$c.Value = "'16.56"
$c.Style = $s
$ws.Range("E11").Value = "  +3.82%  "
$ws.Range("E12").Value = "  +1.00%  "
$c = $ws.Range("D13")
$s = $c.Style
$c.Value = "'0.6293"
$c.Style = $s
$ws.Range("E13").Value = "  +2.86%  "
$c = $ws.Range("D14")
$s = $c.Style
$c.Value = "'4.665"
$c.Style = $s
$ws.Range("E14").Value = "  +3.43%  "
$c = $ws.Range("D15")
$s = $c.Style
$c.Value = "'80.16"
$c.Style = $s
$ws.Range("E15").Value = "  +3.83%  "
$ws.Range("D16").Value = "28.255.83"
$ws.Range("E16").Value = "  +6.52%  "
$c = $ws.Range("D17")
$s = $c.Style
$c.Value = "'0.9991"
$c.Style = $s
$ws.Range("E17").Value = "  +0.07%  "
$c = $ws.Range("D18")
$s = $c.Style
$c.Value = "'0.9989"
$c.Style = $s
$ws.Range("E18").Value = "  +0.05%  "
$c = $ws.Range("D19")
$s = $c.Style
$c.Value = "'0.000007254"
$c.Style = $s
$ws.Range("E19").Value = "  +0.76%  "
$c = $ws.Range("D20")
$s = $c.Style
$c.Value = "'12.07"
$c.Style = $s
$ws.Range("E20").Value = "  +5.75%  "
$ws.Range("D21").Value = "2.015.09"
$c = $ws.Range("D22")
$s = $c.Style
$c.Value = "'4.567"
$c.Style = $s
$ws.Range("E22").Value = "  +2.17%  "
$c = $ws.Range("D23")
$s = $c.Style
$c.Value = "'8.778"
$c.Style = $s
$ws.Range("E23").Value = "  +2.43%  "
$c = $ws.Range("D24")
$s = $c.Style
$c.Value = "'5.253"
$c.Style = $s
$ws.Range("E24").Value = "  +2.90%  "
$c = $ws.Range("D25")
$s = $c.Style
$c.Value = "'142.04"
$c.Style = $s
$ws.Range("E25").Value = "  +2.80%  "
$c = $ws.Range("D26")
$s = $c.Style
$c.Value = "'15.78"
$c.Style = $s
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("E27").Value = "  +4.83%  "
$c = $ws.Range("D28")
$s = $c.Style
$c.Value = "'109.86"
$c.Style = $s
$ws.Range("E28").Value = "  +3.16%  "
$c = $ws.Range("D29")
$s = $c.Style
$c.Value = "'1.385"
$c.Style = $s
$ws.Range("E29").Value = "  +0.14%  "
$c = $ws.Range("D30")
$s = $c.Style
$c.Value = "'4.168"
$c.Style = $s
$ws.Range("E30").Value = "  +5.85%  "
$c = $ws.Range("D31")
$s = $c.Style
$c.Value = "'0.08279"
$c.Style = $s
$ws.Range("E31").Value = "  +3.59%  "
$c = $ws.Range("D32")
$s = $c.Style
$c.Value = "'3.777"
$c.Style = $s
$ws.Range("E32").Value = "  +2.82%  "
$c = $ws.Range("D33")
$s = $c.Style
$c.Value = "'0.04890"
$c.Style = $s
$ws.Range("E33").Value = "  +9.01%  "
$ws.Range("E34").Value = "  +7.81%  "
$ws.Range("E35").Value = "  +5.76%  "
$c = $ws.Range("D36")
$s = $c.Style
$c.Value = "'2.615"
$c.Style = $s
$ws.Range("E36").Value = "  +0.19%  "
$c = $ws.Range("D37")
$s = $c.Style
$c.Value = "'0.9466"
$c.Style = $s
$ws.Range("E37").Value = "  +0.44%  "
$c = $ws.Range("D38")
$s = $c.Style
$c.Value = "'2.611"
$c.Style = $s
$ws.Range("E38").Value = "  +7.79%  "
$c = $ws.Range("D39")
$s = $c.Style
$c.Value = "'2.066"
$c.Style = $s
$ws.Range("E39").Value = "  +0.88%  "
$c = $ws.Range("D40")
$s = $c.Style
$c.Value = "'5.908"
$c.Style = $s
$ws.Range("E40").Value = "  +5.92%  "
$c = $ws.Range("D41")
$s = $c.Style
$c.Value = "'0.01553"
$c.Style = $s
$ws.Range("E41").Value = "  +2.82%  "
$c = $ws.Range("D42")
$s = $c.Style
$c.Value = "'0.9988"
$c.Style = $s
$ws.Range("E42").Value = "  +0.10%  "
$c = $ws.Range("D43")
$s = $c.Style
$c.Value = "'99.76"
$c.Style = $s
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("E44").Value = "  +3.48%  "
$c = $ws.Range("D45")
$s = $c.Style
$c.Value = "'7.199"
$c.Style = $s
$ws.Range("E45").Value = "  +4.22%  "
$c = $ws.Range("D46")
$s = $c.Style
$c.Value = "'0.1213"
$c.Style = $s
$ws.Range("E46").Value = "  +4.48%  "
$c = $ws.Range("D47")
$s = $c.Style
$c.Value = "'0.05451"
$c.Style = $s
$ws.Range("E47").Value = "  +1.31%  "
$c = $ws.Range("D48")
$s = $c.Style
$c.Value = "'8.031"
$c.Style = $s
$ws.Range("E48").Value = "  +1.68%  "
$c = $ws.Range("D49")
$s = $c.Style
$c.Value = "'30.77"
$c.Style = $s
$ws.Range("E49").Value = "  +1.53%  "
$c = $ws.Range("D50")
$s = $c.Style
$c.Value = "'1.295"
$c.Style = $s
$ws.Range("E50").Value = "  +5.00%  "
$c = $ws.Range("D51")
$s = $c.Style
$c.Value = "'53.01"
$c.Style = $s
$ws.Range("E51").Value = "  +2.50%  "
